$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Fill in "OK" status for UC05 (row 9)
$ws.Range("D9").Value = "OK"

# Fill in the date for UC05 (row 9) - keeps existing date style (numFmtId 14)
$ws.Range("E9").Value = (Get-Date -Year 2010 -Month 4 -Day 14 -Hour 0 -Minute 0 -Second 0)

# Update the active selection to F9, matching the saved selection state
$ws.Range("F9").Select()
